$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Insert a new row before the "Tour" section heading to hold the new
# "UpdateUserActiveStatus" team endpoint, keeping the existing blank
# separator row pattern used throughout the sheet.
$ws.Rows.Item(57).Insert()

$ws.Cells.Item(57, 2).Value = "UpdateUserActiveStatus"
$ws.Cells.Item(57, 4).Value = "PUT"

[void]$ws.Range("F35").Select()
